# Apply the "Add files via upload" planning.xlsx edit:
#  - Fill in the 9 "Page_*" rows (model/interface subtasks) in rows 63-80
#  - Fill in the JavaDoc row (81)
#  - Clear the stray formatted-but-empty rows 83 and 88 down to just G/H
#  - Append 7 new blank (but styled) rows 89-95
#  - Update sheet view selection, dimension, and column A width
#  - Extend conditional-formatting / data-validation ranges to cover the
#    newly added rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taches")

# ---------------------------------------------------------------------
# 1) Cell values - written in the same order the strings were first
#    introduced so the shared-string table comes out in the expected
#    order (75..87).
# ---------------------------------------------------------------------

$ws.Range("A63").Value = "Page_Accueil"
$ws.Range("A65").Value = "Page_Plante"
$ws.Range("A81").Value = "JavaDoc"
$ws.Range("B81").Value = "généré la JavaDoc"
$ws.Range("A67").Value = "Page_Glossaire"
$ws.Range("A69").Value = "Page_ChoixParcelle"
$ws.Range("A71").Value = "Page_CréationParcelle"
$ws.Range("A73").Value = "Page_ChoixModePlantation"
$ws.Range("A75").Value = "Page_PlantationRapide"
$ws.Range("A79").Value = "Page_Météo"
$ws.Range("B63").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B64").Value = "réaliser l'affichage (interface)"
$ws.Range("A77").Value = "Page_PlantationAssisté"

# Fill in the rest of the "modèle / interface" row pairs (B/D/E columns);
# A-column labels for the even "interface" rows stay blank, matching the
# existing pattern used by every other page above.
$ws.Range("B65").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B66").Value = "réaliser l'affichage (interface)"
$ws.Range("B67").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B68").Value = "réaliser l'affichage (interface)"
$ws.Range("B69").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B70").Value = "réaliser l'affichage (interface)"
$ws.Range("B71").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B72").Value = "réaliser l'affichage (interface)"
$ws.Range("B73").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B74").Value = "réaliser l'affichage (interface)"
$ws.Range("B75").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B76").Value = "réaliser l'affichage (interface)"
$ws.Range("B77").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B78").Value = "réaliser l'affichage (interface)"
$ws.Range("B79").Value = "faire le lien avec le contenu (le modèle)"
$ws.Range("B80").Value = "réaliser l'affichage (interface)"

# Type (D), Difficulty (E) and Affectation (F) columns.
$ws.Range("D63").Value = "Classe"
$ws.Range("E63").Value = 5
$ws.Range("D64").Value = "Fonction"
$ws.Range("E64").Value = 21

$ws.Range("D65").Value = "Classe"
$ws.Range("E65").Value = 8
$ws.Range("D66").Value = "Fonction"
$ws.Range("E66").Value = 13

$ws.Range("D67").Value = "Classe"
$ws.Range("E67").Value = 13
$ws.Range("D68").Value = "Fonction"
$ws.Range("E68").Value = 21

$ws.Range("D69").Value = "Classe"
$ws.Range("E69").Value = 13
$ws.Range("D70").Value = "Fonction"
$ws.Range("E70").Value = 21

$ws.Range("D71").Value = "Classe"
$ws.Range("E71").Value = 5
$ws.Range("D72").Value = "Fonction"
$ws.Range("E72").Value = 8

$ws.Range("D73").Value = "Classe"
$ws.Range("E73").Value = 5
$ws.Range("D74").Value = "Fonction"
$ws.Range("E74").Value = 5

$ws.Range("D75").Value = "Classe"
$ws.Range("E75").Value = 5
$ws.Range("D76").Value = "Fonction"
$ws.Range("E76").Value = 8

$ws.Range("D77").Value = "Classe"
$ws.Range("E77").Value = 13
$ws.Range("D78").Value = "Fonction"
$ws.Range("E78").Value = 8

$ws.Range("D79").Value = "Classe"
$ws.Range("E79").Value = 13
$ws.Range("D80").Value = "Fonction"
$ws.Range("E80").Value = 13

$ws.Range("D81").Value = "Autre3"
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = "TOUS"

# ---------------------------------------------------------------------
# 2) Rows 83 and 88 lose their empty (but styled) A:F cells - only G/H
#    remain.
# ---------------------------------------------------------------------
$ws.Range("A83:F83").Clear()
$ws.Range("A88:F88").Clear()

# ---------------------------------------------------------------------
# 3) Append 7 new blank rows (89-95) with the same styling as the other
#    blank rows (copy formats down from row 82).
# ---------------------------------------------------------------------
$ws.Range("A82:H82").Copy()
$ws.Range("A89:H95").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Conditional formatting / data validation ranges grow to match the
#    new used range.
# ---------------------------------------------------------------------
$gCf = $ws.Range("G5:G88").FormatConditions.Item(1)
$gCf.ModifyAppliesToRange($ws.Range("G5:G95"))

$eCf = $ws.Range("E5:E88").FormatConditions.Item(1)
$eCf.ModifyAppliesToRange($ws.Range("E5:E95"))

# ---------------------------------------------------------------------
# 5) Sheet view: drop the scrolled topLeftCell and move the selection.
# ---------------------------------------------------------------------
$ws.Range("G85").Select()

# ---------------------------------------------------------------------
# 6) Column A widens slightly (best-fit) to accommodate the new labels.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 23

$wb.Save()
